$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.783.48"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "2.834.87"
$ws.Range("E3").Value = "  +2.21%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'351.54"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "'113.05"
$ws.Range("E6").Value = "  +4.83%  "

$ws.Range("E7").Value = "  +1.74%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'0.623"
$ws.Range("E9").Value = "  +6.42%  "

$ws.Range("D10").Value = "'40.19"
$ws.Range("E10").Value = "  +1.33%  "

$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("D12").Value = "'0.0851"
$ws.Range("E12").Value = "  +2.29%  "

$ws.Range("D13").Value = "'20.03"
$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("E14").Value = "  +3.55%  "

$ws.Range("D15").Value = "3.283.49"
$ws.Range("E15").Value = "  +2.32%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.849.28"
$ws.Range("E16").Value = "  +2.69%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.977"
$ws.Range("E17").Value = "  +5.91%  "

$ws.Range("D18").Value = "51.888.53"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("E19").Value = "  +12.33%  "

$ws.Range("D20").Value = "'7.61"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").Value = "'13.35"
$ws.Range("E21").Value = "  +1.63%  "

$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("E22").Value = "  +1.37%  "

$ws.Range("D23").Value = "'70.54"
$ws.Range("E23").Value = "  +1.10%  "

$ws.Range("D24").Value = "'268.80"
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("E25").Value = "  +2.37%  "

$ws.Range("D26").Value = "'26.31"
$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  +0.86%  "

$ws.Range("D29").Value = "'10.59"
$ws.Range("E29").Value = "  +3.80%  "

$ws.Range("D30").Value = "'38.92"
$ws.Range("E30").Value = "  +6.39%  "

$ws.Range("E31").Value = "  +2.72%  "

$ws.Range("D32").Value = "'6.30"
$ws.Range("E32").Value = "  +1.89%  "

$ws.Range("D33").Value = "'52.80"
$ws.Range("E33").Value = "  +1.76%  "

$ws.Range("D34").Value = "'0.0455"
$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").Value = "'0.0900"
$ws.Range("E35").Value = "  +9.09%  "

$ws.Range("D36").Value = "'5.65"
$ws.Range("E36").Value = "  +2.74%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").Value = "'19.01"
$ws.Range("E38").Value = "  +4.38%  "

$ws.Range("E39").Value = "  +2.84%  "

$ws.Range("E40").Value = "  +3.25%  "

$ws.Range("E41").Value = "  +2.03%  "

$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("D43").Value = "'122.02"
$ws.Range("E43").Value = "  +0.95%  "

$ws.Range("E44").Value = "  +1.34%  "

$ws.Range("E45").Value = "  +1.18%  "

$ws.Range("D46").Value = "2.186.20"
$ws.Range("E46").Value = "  +4.38%  "

$ws.Range("D47").Value = "'3.50"
$ws.Range("E47").Value = "  +8.11%  "

$ws.Range("D48").Value = "'2.50"
$ws.Range("E48").Value = "  +7.87%  "

$ws.Range("E49").Value = "  +21.36%  "

$ws.Range("D50").Value = "'0.967"
$ws.Range("E50").Value = "  +7.56%  "

$ws.Range("D51").Value = "'5.52"
$ws.Range("E51").Value = "  +2.08%  "
